$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 328.25
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H80").Value = 3080171.8
$ws.Range("I80").Value = 335
$ws.Range("J80").Value = 5500043.5
$ws.Range("K80").Value = 1005
$ws.Range("L80").Value = 16500130.5
$ws.Range("M80").Value = -7
$ws.Range("N80").Value = -16502126.5
$ws.Range("H83").Value = 3080171.8
$ws.Range("I83").Value = 335
$ws.Range("J83").Value = 5500043.5
$ws.Range("K83").Value = 3015
$ws.Range("L83").Value = 49500391.5
$ws.Range("M83").Value = 1977
$ws.Range("N83").Value = -49510375.5
$ws.Range("H116").Value = 6988.25
$ws.Range("I116").Value = 4000
$ws.Range("K116").Value = 4000
$ws.Range("M116").Value = -558
$ws.Range("H129").Value = 118444.89
$ws.Range("J129").Value = 130700.17
$ws.Range("L129").Value = 392100.51
$ws.Range("N129").Value = -402100.51

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3152.6924
$ws.Range("I45").Value = 2597.5417
$ws.Range("K45").Value = 2597.5417
$ws.Range("M45").Value = -2220.5417
$ws.Range("H97").Value = 3092
$ws.Range("I97").Value = 3140
$ws.Range("J97").Value = 2900
$ws.Range("K97").Value = 3140
$ws.Range("L97").Value = 2900
$ws.Range("M97").Value = -2644
$ws.Range("N97").Value = -3892
$ws.Range("H101").Value = 30700.5
$ws.Range("J101").Value = 30700.5
$ws.Range("L101").Value = 30700.5
$ws.Range("N101").Value = -37190.5
$ws.Range("H102").Value = 2000
$ws.Range("I102").Value = 2000
$ws.Range("K102").Value = 2000
$ws.Range("M102").Value = -378
$ws.Range("H132").Value = 15304.5
$ws.Range("I132").Value = 1789.1
$ws.Range("K132").Value = 5367.299999999999
$ws.Range("M132").Value = -2837.299999999999
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1450.6757
$ws.Range("I86").Value = 1362.8788
$ws.Range("J86").Value = 2175
$ws.Range("K86").Value = 1362.8788
$ws.Range("L86").Value = 2175
$ws.Range("M86").Value = -239.8788
$ws.Range("N86").Value = -4421
$ws.Range("H89").Value = 1450.6757
$ws.Range("I89").Value = 1362.8788
$ws.Range("J89").Value = 2175
$ws.Range("K89").Value = 6814.394
$ws.Range("L89").Value = 10875
$ws.Range("M89").Value = -1198.394
$ws.Range("N89").Value = -22107
$ws.Range("H94").Value = 1509.2
$ws.Range("I94").Value = 1202.1538
$ws.Range("K94").Value = 1202.1538
$ws.Range("M94").Value = -751.1538

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3436.375
$ws.Range("J62").Value = 3319.4
$ws.Range("L62").Value = 3319.4
$ws.Range("N62").Value = -4567.4
$ws.Range("H65").Value = 3436.375
$ws.Range("J65").Value = 3319.4
$ws.Range("L65").Value = 16597
$ws.Range("N65").Value = -22837
$ws.Range("H122").Value = 2500.375
$ws.Range("I122").Value = 2500.375
$ws.Range("K122").Value = 7501.125
$ws.Range("M122").Value = -5051.125
$ws.Range("H132").Value = 4502.4
$ws.Range("I132").Value = 1582.4
$ws.Range("J132").Value = 7422.4
$ws.Range("K132").Value = 4747.200000000001
$ws.Range("L132").Value = 22267.2
$ws.Range("M132").Value = -2217.200000000001
$ws.Range("N132").Value = -27327.2
$ws.Range("H134").Value = 1175.7059
$ws.Range("I134").Value = 915.5833
$ws.Range("J134").Value = 1800
$ws.Range("K134").Value = 2746.7499
$ws.Range("L134").Value = 5400
$ws.Range("M134").Value = -211.7498999999998
$ws.Range("N134").Value = -10470

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 41.23077
$ws.Range("I2").Value = 32.47059
$ws.Range("J2").Value = 57.77778
$ws.Range("K2").Value = 194.82354
$ws.Range("L2").Value = 346.66668
$ws.Range("M2").Value = -81.82354000000001
$ws.Range("N2").Value = -572.66668
$ws.Range("H5").Value = 1756
$ws.Range("J5").Value = 2000
$ws.Range("L5").Value = 6000
$ws.Range("N5").Value = -6224
$ws.Range("H107").Value = 33666.668
$ws.Range("H122").Value = 448
$ws.Range("I122").Value = 264.2857
$ws.Range("J122").Value = 769.5
$ws.Range("K122").Value = 2378.5713
$ws.Range("L122").Value = 6925.5
$ws.Range("M122").Value = 71.42869999999994
$ws.Range("N122").Value = -11825.5
$ws.Range("H130").Value = 2999
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 2999
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 8997
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -19037
$ws.Range("H131").Value = 780.47
$ws.Range("J131").Value = 781.28284
$ws.Range("L131").Value = 2343.84852
$ws.Range("N131").Value = -12423.84852
$ws.Range("H135").Value = 1756
$ws.Range("J135").Value = 2000
$ws.Range("L135").Value = 18000
$ws.Range("N135").Value = -23070

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()
$ws.Range("H102").Value = 1894.6923
$ws.Range("I102").Value = 1902.6957
$ws.Range("J102").Value = 1833.3334
$ws.Range("K102").Value = 1902.6957
$ws.Range("L102").Value = 1833.3334
$ws.Range("M102").Value = -280.6957
$ws.Range("N102").Value = -5077.3334
$ws.Range("H132").Value = 46266.5
$ws.Range("I132").Value = 5025
$ws.Range("J132").Value = 128749.5
$ws.Range("K132").Value = 15075
$ws.Range("L132").Value = 386248.5
$ws.Range("M132").Value = -12545
$ws.Range("N132").Value = -391308.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 257.41666
$ws.Range("I55").Value = 175.57143
$ws.Range("J55").Value = 372
$ws.Range("K55").Value = 175.57143
$ws.Range("L55").Value = 372
$ws.Range("M55").Value = -2.571429999999992
$ws.Range("N55").Value = -718
$ws.Range("H68").Value = 2816.7273
$ws.Range("I68").Value = 2837.6
$ws.Range("J68").Value = 2799.3333
$ws.Range("K68").Value = 2837.6
$ws.Range("L68").Value = 2799.3333
$ws.Range("M68").Value = -2088.6
$ws.Range("N68").Value = -4297.3333
$ws.Range("H71").Value = 2816.7273
$ws.Range("I71").Value = 2837.6
$ws.Range("J71").Value = 2799.3333
$ws.Range("K71").Value = 14188
$ws.Range("L71").Value = 13996.6665
$ws.Range("M71").Value = -10444
$ws.Range("N71").Value = -21484.6665

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1782.1177
$ws.Range("I122").Value = 1772
$ws.Range("K122").Value = 5316
$ws.Range("M122").Value = -2866
$ws.Range("H126").Value = 1230.0869
$ws.Range("I126").Value = 717.4545000000001
$ws.Range("K126").Value = 2152.3635
$ws.Range("M126").Value = 317.6364999999996
